$d = $word.ActiveDocument

# --- Edit 1: split the first paragraph's sentence and append a red "(This is a
# change ... )" comment split across three runs with color C00000 ---
$p1 = $d.Paragraphs.Item(1)
$p1Range = $p1.Range
# Range covering just the sentence text (exclude the trailing paragraph mark)
$sentenceRange = $d.Range($p1Range.Start, $p1Range.End - 1)
$sentenceRange.InsertAfter("  ")

$run1 = $d.Range($sentenceRange.End, $sentenceRange.End)
$run1.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$run1.Font.Color = 192

$run2 = $d.Range($run1.End, $run1.End)
$run2.InsertAfter("rsion for branch alternate")
$run2.Font.Color = 192

$run3 = $d.Range($run2.End, $run2.End)
$run3.InsertAfter(")")
$run3.Font.Color = 192

# --- Edit 2: append a new, empty, shaded paragraph after the final paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$endRange.InsertXML($newParaXml)
